$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 12), copying row 2's formatting (thin border,
# vertical-centered) down onto the new row so it matches the rest of the table.
$ws.Range("A2:C2").Copy() | Out-Null
$ws.Range("A12:C12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A12").Value = "new_item1"
$ws.Range("B12").Value = "new_item2"
$ws.Range("C12").Value = "new_item3"

# Set explicit column widths for A:C (matches width="11", customWidth="1")
$ws.Columns("A:C").ColumnWidth = 10.2857142857143

# Update the active selection to match the recorded end state
$ws.Range("C22").Select()
